$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill_pattern_group")

# Copy formatting + values from the first data row (same D/F variants as the
# new "elite warrior" rows), then overwrite with the new row's own data.
$ws.Range("A4:G4").Copy($ws.Range("A17:G17"))
$ws.Range("A4:G4").Copy($ws.Range("A18:G18"))

$ws.Range("A17").Value = 2010101
$ws.Range("B17").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C17").Value = 1
$ws.Range("G17").Value = 201010101

$ws.Range("A18").Value = 2010102
$ws.Range("B18").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C18").Value = 2
$ws.Range("G18").Value = 201010201

$ws.Range("K15").Select()
